$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "_GoBack" bookmark and the
# "Add the SSH key to your GitHub account" hyperlink (paragraph 10).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Add the SSH key to your GitHub account*") {
        $target = $cand
        break
    }
}

$r = $target.Range

# Rewrite the paragraph: drop the (hidden) _GoBack bookmark pair that used to
# sit right after the paragraph properties, and give the paragraph mark its
# own run-properties (bold, blue, SimSun) instead.
$xmlPara = @'
<?xml version="1.0" encoding="utf-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="7"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="SimSun" w:hAnsi="SimSun" w:eastAsia="SimSun" w:cs="SimSun"/><w:b/><w:color w:val="0070C0"/><w:kern w:val="0"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="SimSun" w:hAnsi="SimSun" w:eastAsia="SimSun" w:cs="SimSun"/><w:b/><w:color w:val="0070C0"/><w:kern w:val="0"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:eastAsia="zh-CN" w:bidi="ar"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="SimSun" w:hAnsi="SimSun" w:eastAsia="SimSun" w:cs="SimSun"/><w:b/><w:color w:val="0070C0"/><w:kern w:val="0"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="SimSun" w:hAnsi="SimSun" w:eastAsia="SimSun" w:cs="SimSun"/><w:b/><w:color w:val="0070C0"/><w:kern w:val="0"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK &quot;https://help.github.com/articles/adding-a-new-ssh-key-to-your-github-account&quot; </w:instrText></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="SimSun" w:hAnsi="SimSun" w:eastAsia="SimSun" w:cs="SimSun"/><w:b/><w:color w:val="0070C0"/><w:kern w:val="0"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="SimSun" w:hAnsi="SimSun" w:eastAsia="SimSun" w:cs="SimSun"/><w:b/><w:color w:val="0070C0"/><w:kern w:val="0"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar"/></w:rPr><w:t>Add the SSH key to your GitHub account</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia" w:ascii="SimSun" w:hAnsi="SimSun" w:eastAsia="SimSun" w:cs="SimSun"/><w:b/><w:color w:val="0070C0"/><w:kern w:val="0"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:lang w:val="en-US" w:eastAsia="zh-CN" w:bidi="ar"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xmlPara)

# Re-fetch the (now bookmark-free) paragraph and append two brand-new
# paragraphs after it: "Push" and the "$ git push origin master" command
# line, the latter carrying the relocated _GoBack bookmark at its end.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Add the SSH key to your GitHub account*") {
        $target = $cand
        break
    }
}

$tail = $target.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$pushPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -eq $tail.End) {
        $pushPara = $cand
        break
    }
}

$xmlPush = @'
<?xml version="1.0" encoding="utf-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="7"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:ind w:left="0" w:leftChars="0" w:firstLine="0" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorAscii" w:hAnsiTheme="minorAscii" w:eastAsiaTheme="minorEastAsia" w:cstheme="minorBidi"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:asciiTheme="minorAscii" w:hAnsiTheme="minorAscii" w:eastAsiaTheme="minorEastAsia" w:cstheme="minorBidi"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>Push</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$pushPara.Range.InsertXML($xmlPush)

# Re-fetch the "Push" paragraph and add the git-push command paragraph
# right after it, ending with the relocated _GoBack bookmark.
$pushPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Push`r") {
        $pushPara = $cand
    }
}

$tail2 = $pushPara.Range
$tail2.Collapse(0)
$tail2.InsertParagraphAfter()

$cmdPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -eq $tail2.End) {
        $cmdPara = $cand
        break
    }
}

$xmlCmd = @'
<?xml version="1.0" encoding="utf-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="7"/><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:suppressLineNumbers w:val="0"/><w:ind w:left="420" w:leftChars="0" w:firstLine="420" w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia" w:asciiTheme="minorAscii" w:hAnsiTheme="minorAscii" w:eastAsiaTheme="minorEastAsia" w:cstheme="minorBidi"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:cstheme="minorBidi"/><w:kern w:val="2"/><w:sz w:val="21"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN" w:bidi="ar-SA"/></w:rPr><w:t>$ git push origin master</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$cmdPara.Range.InsertXML($xmlCmd)
